$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need NumberFormat "@" (Text)
# applied first, matching the source data which stores these as literal text
# (e.g. trailing zeros / fixed decimal places) rather than numeric values.
$textCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D11", "D12", "D13", "D14", "D15", "D17", "D18", "D19", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '31.004.90'
$ws.Range("E2").Value = '  +3.34%  '
$ws.Range("D3").Value = '1.909.86'
$ws.Range("E3").Value = '  +1.47%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").Value = '245.27'
$ws.Range("E5").Value = '  +0.53%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.23%  '
$ws.Range("D7").Value = '0.4986'
$ws.Range("E7").Value = '  +0.61%  '
$ws.Range("D8").Value = '0.2991'
$ws.Range("E8").Value = '  +2.46%  '
$ws.Range("D9").Value = '0.06864'
$ws.Range("E9").Value = '  +3.74%  '
$ws.Range("D10").Value = '1.914.97'
$ws.Range("E10").Value = '  +1.91%  '
$ws.Range("D11").Value = '16.98'
$ws.Range("E11").Value = '  +0.13%  '
$ws.Range("D12").Value = '0.07302'
$ws.Range("E12").Value = '  +1.52%  '
$ws.Range("D13").Value = '91.14'
$ws.Range("E13").Value = '  +6.11%  '
$ws.Range("D14").Value = '5.098'
$ws.Range("E14").Value = '  +5.20%  '
$ws.Range("D15").Value = '0.6798'
$ws.Range("E15").Value = '  +2.02%  '
$ws.Range("D16").Value = '30.975.08'
$ws.Range("E16").Value = '  +3.30%  '
$ws.Range("D17").Value = '0.000008049'
$ws.Range("E17").Value = '  +2.93%  '
$ws.Range("B18").Value = 'Avalanche'
$ws.Range("C18").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D18").Value = '13.29'
$ws.Range("E18").Value = '  +3.67%  '
$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D19").Value = '1.001'
$ws.Range("E19").Value = '  +0.36%  '
$ws.Range("D20").Value = '2.161.60'
$ws.Range("E20").Value = '  +1.89%  '
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  +0.21%  '
$ws.Range("D22").Value = '4.880'
$ws.Range("E22").Value = '  +2.41%  '
$ws.Range("D23").Value = '184.13'
$ws.Range("E23").Value = '  +35.58%  '
$ws.Range("D24").Value = '6.108'
$ws.Range("E24").Value = '  +8.97%  '
$ws.Range("D25").Value = '9.369'
$ws.Range("E25").Value = '  +2.24%  '
$ws.Range("D26").Value = '153.90'
$ws.Range("E26").Value = '  +2.82%  '
$ws.Range("D27").Value = '18.74'
$ws.Range("E27").Value = '  +11.85%  '
$ws.Range("D28").Value = '1.941'
$ws.Range("E28").Value = '  +1.82%  '
$ws.Range("D29").Value = '1.401'
$ws.Range("E29").Value = '  +1.40%  '
$ws.Range("D30").Value = '4.349'
$ws.Range("E30").Value = '  +4.07%  '
$ws.Range("D31").Value = '0.08971'
$ws.Range("E31").Value = '  +3.56%  '
$ws.Range("D32").Value = '4.048'
$ws.Range("E32").Value = '  +2.38%  '
$ws.Range("D33").Value = '0.05251'
$ws.Range("E33").Value = '  +5.35%  '
$ws.Range("D34").Value = '0.7494'
$ws.Range("E34").Value = '  +6.70%  '
$ws.Range("D35").Value = '1.142'
$ws.Range("E35").Value = '  +3.37%  '
$ws.Range("D36").Value = '2.669'
$ws.Range("E36").Value = '  +0.60%  '
$ws.Range("D37").Value = '0.01939'
$ws.Range("E37").Value = '  +18.15%  '
$ws.Range("D38").Value = '2.737'
$ws.Range("E38").Value = '  +1.64%  '
$ws.Range("D39").Value = '2.185'
$ws.Range("E39").Value = '  -0.30%  '
$ws.Range("D40").Value = '0.9367'
$ws.Range("E40").Value = '  +0.56%  '
$ws.Range("D41").Value = '0.4381'
$ws.Range("E41").Value = '  +4.83%  '
$ws.Range("D42").Value = '106.11'
$ws.Range("E42").Value = '  +4.81%  '
$ws.Range("D43").Value = '5.881'
$ws.Range("E43").Value = '  -1.41%  '
$ws.Range("D44").Value = '1.001'
$ws.Range("E44").Value = '  +0.14%  '
$ws.Range("D45").Value = '7.802'
$ws.Range("E45").Value = '  +3.40%  '
$ws.Range("D46").Value = '0.1346'
$ws.Range("E46").Value = '  +6.91%  '
$ws.Range("D47").Value = '0.05857'
$ws.Range("E47").Value = '  +2.26%  '
$ws.Range("D48").Value = '8.584'
$ws.Range("E48").Value = '  +4.44%  '
$ws.Range("D49").Value = '0.3899'
$ws.Range("E49").Value = '  +5.35%  '
$ws.Range("D50").Value = '33.28'
$ws.Range("E50").Value = '  +2.95%  '
$ws.Range("E51").Value = '  +4.29%  '
